# Update gh-pages output (苏州-漫展信息.xlsx) to data generated at 456a3b4.
# Sheet 1 "展览" and Sheet 4 "全部类型" both carry the same event table
# (sheet 4 has one extra inserted row, so its row numbers are offset by +1
# from row 33 onward). The refreshed scrape bumped "interested" counts
# (column F, "想去人数") for a batch of rows, and for row 3 specifically
# the minimum price (column G, "最低票价") went from 89 to 99.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("G3").Value = 99

$ws1.Range("F4").Value = 1151
$ws1.Range("F5").Value = 131
$ws1.Range("F6").Value = 104
$ws1.Range("F9").Value = 1189
$ws1.Range("F10").Value = 16848
$ws1.Range("F14").Value = 6494
$ws1.Range("F15").Value = 659
$ws1.Range("F19").Value = 131
$ws1.Range("F20").Value = 1283
$ws1.Range("F21").Value = 95
$ws1.Range("F22").Value = 42
$ws1.Range("F23").Value = 645
$ws1.Range("F24").Value = 26
$ws1.Range("F28").Value = 918
$ws1.Range("F29").Value = 77
$ws1.Range("F30").Value = 5082
$ws1.Range("F33").Value = 11533
$ws1.Range("F34").Value = 1259
$ws1.Range("F37").Value = 234
$ws1.Range("F39").Value = 277
$ws1.Range("F40").Value = 79

# --- Sheet 4: 全部类型 (same rows, offset by +1 from row 33 down) ----
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("G3").Value = 99

$ws4.Range("F4").Value = 1151
$ws4.Range("F5").Value = 131
$ws4.Range("F6").Value = 104
$ws4.Range("F9").Value = 1189
$ws4.Range("F10").Value = 16848
$ws4.Range("F14").Value = 6494
$ws4.Range("F15").Value = 659
$ws4.Range("F19").Value = 131
$ws4.Range("F20").Value = 1283
$ws4.Range("F21").Value = 95
$ws4.Range("F22").Value = 42
$ws4.Range("F23").Value = 645
$ws4.Range("F24").Value = 26
$ws4.Range("F28").Value = 918
$ws4.Range("F29").Value = 77
$ws4.Range("F30").Value = 5082
$ws4.Range("F34").Value = 11533
$ws4.Range("F35").Value = 1259
$ws4.Range("F38").Value = 234
$ws4.Range("F40").Value = 277
$ws4.Range("F41").Value = 79
